# Keyword_Frequency__zh-en-comparison.xlsx — apply manual-comparison edits
# (compared to recovery PDF from 2020-09-27 10:00am)

$wb = $excel.ActiveWorkbook

$wsPos = $wb.Worksheets.Item("Positive")
$wsNeg = $wb.Worksheets.Item("Negative")

# --- Negative sheet: fill in previously-blank Chinese-keyword / count cells ---
# Row 11 (5: 15,000 to 20,000 yen block)
$wsNeg.Range("B11").Value = "\begin{CJK}{UTF8}{gbsn}华人\end{CJK} (Chinese)"
$wsNeg.Range("C11").Value = 15

# Row 19 (6: 20,000 to 30,000 yen block)
$wsNeg.Range("B19").Value = "\begin{CJK}{UTF8}{gbsn}老\end{CJK} (old)"
$wsNeg.Range("C19").Value = 2

# Row 41 (8: 50,000 to 100,000 yen block)
$wsNeg.Range("B41").Value = "\begin{CJK}{UTF8}{gbsn}华人\end{CJK} (Chinese)"
$wsNeg.Range("C41").Value = 3

# Row 51 (9: 100,000 to 200,000 yen block)
$wsNeg.Range("B51").Value = "\begin{CJK}{UTF8}{gbsn}华人\end{CJK} (Chinese)"
$wsNeg.Range("C51").Value = 2

# Row 61
$wsNeg.Range("B61").Value = "\begin{CJK}{UTF8}{gbsn}华人\end{CJK} (Chinese)"
$wsNeg.Range("C61").Value = 2

# Row 71
$wsNeg.Range("B71").Value = "\begin{CJK}{UTF8}{gbsn}华人\end{CJK} (Chinese)"
$wsNeg.Range("C71").Value = 8

# --- Positive sheet: header row (row 3) loses its centered horizontal
#     alignment (kept vertical centering where it previously existed) ---
$wsPos.Range("A3:D3").HorizontalAlignment = 1   # xlHAlignGeneral
$wsPos.Range("E3").HorizontalAlignment = 1      # xlHAlignGeneral (keeps existing vertical="center")

# --- Window / selection state ---
# Negative tab was active/selected before; now Positive is active and
# each sheet remembers a different selected cell.
$wsNeg.Range("H16").Select()
$wsPos.Activate()
$wsPos.Range("G22").Select()
